# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) Slide Master, clrScheme "Red
#                             Violet" (theme display name "Integral")
#   ppt/theme/theme2.xml  -> bound to the Notes Master, clrScheme "Office"
#                             (theme display name "Office Theme")
# The target revision swaps the two clrScheme colour palettes between the
# parts (font/format schemes are already identical between the two themes,
# so only the 12 theme colours actually change). Concretely theme1.xml ends
# up holding the plain "Office" palette.
#
# PowerPoint's ColorScheme.Colors(index) exposes exactly those 12 theme
# colours, in this fixed order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3  8 accent4
#   9 accent5  10 accent6  11 hlink  12 folHlink
# RGBColor.RGB takes a COM "RGB()" style integer, i.e. R + G*256 + B*65536
# (the bytes of a VBA RGB() call), not a straightforward hex RRGGBB number.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0         # dk1      000000
$scheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388   # dk2      44546A
$scheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407     # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$scheme.Colors(10).RGB = 4697456   # accent6  70AD47
$scheme.Colors(11).RGB = 12673797  # hlink    0563C1
$scheme.Colors(12).RGB = 7491477   # folHlink 954F72
